# Insert two new weekly rows of data for "Albahaca" right after row 534,
# shifting the existing rows 535-545 down to 537-547.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at position 535 (pushes old 535:545 down to 537:547)
$ws.Range("A535:A536").EntireRow.Insert()

# --- New row 535 ---
$ws.Cells.Item(535, 1).Value = 6
$ws.Cells.Item(535, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(535, 3).Value = "Metropolitana"
$ws.Cells.Item(535, 4).Value = 44890
$ws.Cells.Item(535, 5).Value = 13
$ws.Cells.Item(535, 6).Value = 100112052
$ws.Cells.Item(535, 7).Value = "Albahaca"
$ws.Cells.Item(535, 8).Value = "Sin especificar"
$ws.Cells.Item(535, 9).Value = "Primera"
$ws.Cells.Item(535, 10).Value = 580
$ws.Cells.Item(535, 11).Value = 4500
$ws.Cells.Item(535, 12).Value = 5000
$ws.Cells.Item(535, 13).Value = 4776
$ws.Cells.Item(535, 14).Value = "`$/docena de matas"
$ws.Cells.Item(535, 15).Value = "Región Metropolitana"
$ws.Cells.Item(535, 16).Value = 796
$ws.Cells.Item(535, 17).Value = 6
$ws.Cells.Item(535, 18).Value = "Hortaliza"

# --- New row 536 ---
$ws.Cells.Item(536, 1).Value = 6
$ws.Cells.Item(536, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(536, 3).Value = "Metropolitana"
$ws.Cells.Item(536, 4).Value = 44890
$ws.Cells.Item(536, 5).Value = 13
$ws.Cells.Item(536, 6).Value = 100112052
$ws.Cells.Item(536, 7).Value = "Albahaca"
$ws.Cells.Item(536, 8).Value = "Sin especificar"
$ws.Cells.Item(536, 9).Value = "Segunda"
$ws.Cells.Item(536, 10).Value = 150
$ws.Cells.Item(536, 11).Value = 4000
$ws.Cells.Item(536, 12).Value = 4000
$ws.Cells.Item(536, 13).Value = 4000
$ws.Cells.Item(536, 14).Value = "`$/docena de matas"
$ws.Cells.Item(536, 15).Value = "Región Metropolitana"
$ws.Cells.Item(536, 16).Value = 667
$ws.Cells.Item(536, 17).Value = 6
$ws.Cells.Item(536, 18).Value = "Hortaliza"

Write-Host ("UsedRange: " + $ws.UsedRange.Address())
